$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update header values for columns B:E
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: clear B2, D2, E2 entirely (deleted values); update C2
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -0.5911246880189821
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3: update B3:E3 values
$ws.Range("B3").Value = -0.75226107008933984
$ws.Range("C3").Value = -0.026578298959722881
$ws.Range("D3").Value = -1.5154658291482421
$ws.Range("E3").Value = 2.1862208822635893

# Update the selection to match the new sqref B1:E3
$ws.Range("B1:E3").Select()
